# Updating graphs. Solved some smaller formatting issues.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("slide.pot.objects")

# Column F values: 7.76 -> 9.1999999999999993 (effectively 9.2) across the
# rows that held the old constant.
$fRows = @(2, 3, 4, 5, 16, 18, 19, 20, 23, 24, 26, 27, 35, 37, 39)
foreach ($r in $fRows) {
    $ws.Cells.Item($r, 6).Value = 9.1999999999999993
}

# J7: 20 -> 14
$ws.Cells.Item(7, 10).Value = 14

# Update the active selection on this sheet from N5 to F14.
$ws.Range("F14").Select()

$wb.Save()
